$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.001628901849080777
$ws.Range("J2").Value = 0.001628901849080777
$ws.Range("M2").Value = 4.824089
$ws.Range("N2").Value = 14.472267
$ws.Range("O2").Value = 0.1181976021471384
$ws.Range("P2").Value = 0.1181976021471384
$ws.Range("Q2").Value = 0.4517550304643334
$ws.Range("R2").Value = 4.065795274179
$ws.Range("S2").Value = 0.0001925322926943878
$ws.Range("T2").Value = 0.0001925322926943878
$ws.Range("I3").Value = 0.001628901849080777
$ws.Range("J3").Value = 0.001628901849080777
$ws.Range("M3").Value = 0.7217316666666668
$ws.Range("O3").Value = 0.01768353618551768
$ws.Range("P3").Value = 0.01768353618551768
$ws.Range("Q3").Value = 0.06758704307944445
$ws.Range("R3").Value = 0.608283387715
$ws.Range("S3").Value = 0.00002880474479087658
$ws.Range("T3").Value = 0.00002880474479087657
$ws.Range("I4").Value = 0.001628901849080777
$ws.Range("J4").Value = 0.001628901849080777
$ws.Range("M4").Value = 4.718667333333333
$ws.Range("N4").Value = 14.156002
$ws.Range("O4").Value = 0.1156146091272429
$ws.Range("P4").Value = 0.1156146091272429
$ws.Range("Q4").Value = 0.4418827482082222
$ws.Range("R4").Value = 3.976944733874
$ws.Range("S4").Value = 0.0001883248505881172
$ws.Range("T4").Value = 0.0001883248505881172
$ws.Range("I5").Value = 0.001628901849080777
$ws.Range("J5").Value = 0.001628901849080777
$ws.Range("M5").Value = 27.97761966666667
$ws.Range("N5").Value = 83.93285900000001
$ws.Range("O5").Value = 0.6854947241613126
$ws.Range("P5").Value = 0.6854947241613125
$ws.Range("Q5").Value = 2.619982845431445
$ws.Range("R5").Value = 23.579845608883
$ws.Range("S5").Value = 0.001116603623721479
$ws.Range("T5").Value = 0.001116603623721479
$ws.Range("I6").Value = 0.001628901849080777
$ws.Range("J6").Value = 0.001628901849080777
$ws.Range("M6").Value = 2.571656
$ws.Range("N6").Value = 7.714968
$ws.Range("O6").Value = 0.0630095283787885
$ws.Range("P6").Value = 0.06300952837878848
$ws.Range("Q6").Value = 0.2408244405573333
$ws.Range("R6").Value = 2.167419965016
$ws.Range("S6").Value = 0.0001026363372859163
$ws.Range("T6").Value = 0.0001026363372859162
$ws.Range("I7").Value = 0.1785947081647151
$ws.Range("J7").Value = 0.178594708164715
$ws.Range("M7").Value = 4.824089
$ws.Range("N7").Value = 14.472267
$ws.Range("O7").Value = 0.1181976021471384
$ws.Range("P7").Value = 0.1181976021471384
$ws.Range("Q7").Value = 49.53095109643933
$ws.Range("R7").Value = 445.778559867954
$ws.Range("S7").Value = 0.02110946626123728
$ws.Range("T7").Value = 0.02110946626123728
$ws.Range("I8").Value = 0.1785947081647151
$ws.Range("J8").Value = 0.178594708164715
$ws.Range("M8").Value = 0.7217316666666668
$ws.Range("O8").Value = 0.01768353618551768
$ws.Range("P8").Value = 0.01768353618551768
$ws.Range("Q8").Value = 7.410322630121112
$ws.Range("R8").Value = 66.69290367109001
$ws.Range("S8").Value = 0.003158185984372709
$ws.Range("T8").Value = 0.003158185984372708
$ws.Range("I9").Value = 0.1785947081647151
$ws.Range("J9").Value = 0.178594708164715
$ws.Range("M9").Value = 4.718667333333333
$ws.Range("N9").Value = 14.156002
$ws.Range("O9").Value = 0.1156146091272429
$ws.Range("P9").Value = 0.1156146091272429
$ws.Range("Q9").Value = 48.44854249739156
$ws.Range("R9").Value = 436.036882476524
$ws.Range("S9").Value = 0.02064815737665754
$ws.Range("T9").Value = 0.02064815737665754
$ws.Range("I10").Value = 0.1785947081647151
$ws.Range("J10").Value = 0.178594708164715
$ws.Range("M10").Value = 27.97761966666667
$ws.Range("N10").Value = 83.93285900000001
$ws.Range("O10").Value = 0.6854947241613126
$ws.Range("P10").Value = 0.6854947241613125
$ws.Range("Q10").Value = 287.2579903696731
$ws.Range("R10").Value = 2585.321913327058
$ws.Range("S10").Value = 0.1224257302100415
$ws.Range("T10").Value = 0.1224257302100414
$ws.Range("I11").Value = 0.1785947081647151
$ws.Range("J11").Value = 0.178594708164715
$ws.Range("M11").Value = 2.571656
$ws.Range("N11").Value = 7.714968
$ws.Range("O11").Value = 0.0630095283787885
$ws.Range("P11").Value = 0.06300952837878848
$ws.Range("Q11").Value = 26.40427396195733
$ws.Range("R11").Value = 237.638465657616
$ws.Range("S11").Value = 0.01125316833240606
$ws.Range("T11").Value = 0.01125316833240606
$ws.Range("G12").Value = 23.67539566666666
$ws.Range("H12").Value = 71.026187
$ws.Range("I12").Value = 0.4118171950916292
$ws.Range("J12").Value = 0.4118171950916292
$ws.Range("M12").Value = 4.824089
$ws.Range("N12").Value = 14.472267
$ws.Range("O12").Value = 0.1181976021471384
$ws.Range("P12").Value = 0.1181976021471384
$ws.Range("Q12").Value = 114.2122158062143
$ws.Range("R12").Value = 1027.909942255929
$ws.Range("S12").Value = 0.04867580498279087
$ws.Range("T12").Value = 0.04867580498279087
$ws.Range("G13").Value = 23.67539566666666
$ws.Range("H13").Value = 71.026187
$ws.Range("I13").Value = 0.4118171950916292
$ws.Range("J13").Value = 0.4118171950916292
$ws.Range("M13").Value = 0.7217316666666668
$ws.Range("O13").Value = 0.01768353618551768
$ws.Range("P13").Value = 0.01768353618551768
$ws.Range("Q13").Value = 17.08728277349611
$ws.Range("R13").Value = 153.785544961465
$ws.Range("S13").Value = 0.007282384271221219
$ws.Range("T13").Value = 0.007282384271221218
$ws.Range("G14").Value = 23.67539566666666
$ws.Range("H14").Value = 71.026187
$ws.Range("I14").Value = 0.4118171950916292
$ws.Range("J14").Value = 0.4118171950916292
$ws.Range("M14").Value = 4.718667333333333
$ws.Range("N14").Value = 14.156002
$ws.Range("O14").Value = 0.1156146091272429
$ws.Range("P14").Value = 0.1156146091272429
$ws.Range("Q14").Value = 111.7163161360415
$ws.Range("R14").Value = 1005.446845224374
$ws.Range("S14").Value = 0.04761208404239623
$ws.Range("T14").Value = 0.04761208404239622
$ws.Range("G15").Value = 23.67539566666666
$ws.Range("H15").Value = 71.026187
$ws.Range("I15").Value = 0.4118171950916292
$ws.Range("J15").Value = 0.4118171950916292
$ws.Range("M15").Value = 27.97761966666667
$ws.Range("N15").Value = 83.93285900000001
$ws.Range("O15").Value = 0.6854947241613126
$ws.Range("P15").Value = 0.6854947241613125
$ws.Range("Q15").Value = 662.3812154198481
$ws.Range("R15").Value = 5961.430938778633
$ws.Range("S15").Value = 0.2822985145542218
$ws.Range("T15").Value = 0.2822985145542218
$ws.Range("G16").Value = 23.67539566666666
$ws.Range("H16").Value = 71.026187
$ws.Range("I16").Value = 0.4118171950916292
$ws.Range("J16").Value = 0.4118171950916292
$ws.Range("M16").Value = 2.571656
$ws.Range("N16").Value = 7.714968
$ws.Range("O16").Value = 0.0630095283787885
$ws.Range("P16").Value = 0.06300952837878848
$ws.Range("Q16").Value = 60.88497331855732
$ws.Range("R16").Value = 547.9647598670159
$ws.Range("S16").Value = 0.02594840724099909
$ws.Range("T16").Value = 0.02594840724099908
$ws.Range("G17").Value = 0.3314846666666666
$ws.Range("H17").Value = 0.994454
$ws.Range("I17").Value = 0.005765947381177186
$ws.Range("J17").Value = 0.005765947381177185
$ws.Range("M17").Value = 4.824089
$ws.Range("N17").Value = 14.472267
$ws.Range("O17").Value = 0.1181976021471384
$ws.Range("P17").Value = 0.1181976021471384
$ws.Range("Q17").Value = 1.599111534135333
$ws.Range("R17").Value = 14.392003807218
$ws.Range("S17").Value = 0.0006815211545617157
$ws.Range("T17").Value = 0.0006815211545617156
$ws.Range("G18").Value = 0.3314846666666666
$ws.Range("H18").Value = 0.994454
$ws.Range("I18").Value = 0.005765947381177186
$ws.Range("J18").Value = 0.005765947381177185
$ws.Range("M18").Value = 0.7217316666666668
$ws.Range("O18").Value = 0.01768353618551768
$ws.Range("P18").Value = 0.01768353618551768
$ws.Range("Q18").Value = 0.2392429809477778
$ws.Range("R18").Value = 2.15318682853
$ws.Range("S18").Value = 0.0001019623391588377
$ws.Range("T18").Value = 0.0001019623391588376
$ws.Range("G19").Value = 0.3314846666666666
$ws.Range("H19").Value = 0.994454
$ws.Range("I19").Value = 0.005765947381177186
$ws.Range("J19").Value = 0.005765947381177185
$ws.Range("M19").Value = 4.718667333333333
$ws.Range("N19").Value = 14.156002
$ws.Range("O19").Value = 0.1156146091272429
$ws.Range("P19").Value = 0.1156146091272429
$ws.Range("Q19").Value = 1.564165868100889
$ws.Range("R19").Value = 14.077492812908
$ws.Range("S19").Value = 0.00066662775272305
$ws.Range("T19").Value = 0.0006666277527230498
$ws.Range("G20").Value = 0.3314846666666666
$ws.Range("H20").Value = 0.994454
$ws.Range("I20").Value = 0.005765947381177186
$ws.Range("J20").Value = 0.005765947381177185
$ws.Range("M20").Value = 27.97761966666667
$ws.Range("N20").Value = 83.93285900000001
$ws.Range("O20").Value = 0.6854947241613126
$ws.Range("P20").Value = 0.6854947241613125
$ws.Range("Q20").Value = 9.274151929331778
$ws.Range("R20").Value = 83.467367363986
$ws.Range("S20").Value = 0.003952526509588698
$ws.Range("T20").Value = 0.003952526509588696
$ws.Range("G21").Value = 0.3314846666666666
$ws.Range("H21").Value = 0.994454
$ws.Range("I21").Value = 0.005765947381177186
$ws.Range("J21").Value = 0.005765947381177185
$ws.Range("M21").Value = 2.571656
$ws.Range("N21").Value = 7.714968
$ws.Range("O21").Value = 0.0630095283787885
$ws.Range("P21").Value = 0.06300952837878848
$ws.Range("Q21").Value = 0.8524645319413333
$ws.Range("R21").Value = 7.672180787472
$ws.Range("S21").Value = 0.0003633096251448851
$ws.Range("T21").Value = 0.0003633096251448849
$ws.Range("G22").Value = 23.12211433333333
$ws.Range("H22").Value = 69.366343
$ws.Range("I22").Value = 0.4021932475133977
$ws.Range("J22").Value = 0.4021932475133977
$ws.Range("M22").Value = 4.824089
$ws.Range("N22").Value = 14.472267
$ws.Range("O22").Value = 0.1181976021471384
$ws.Range("P22").Value = 0.1181976021471384
$ws.Range("Q22").Value = 111.5431374121757
$ws.Range("R22").Value = 1003.888236709581
$ws.Range("S22").Value = 0.04753827745585415
$ws.Range("T22").Value = 0.04753827745585414
$ws.Range("G23").Value = 23.12211433333333
$ws.Range("H23").Value = 69.366343
$ws.Range("I23").Value = 0.4021932475133977
$ws.Range("J23").Value = 0.4021932475133977
$ws.Range("M23").Value = 0.7217316666666668
$ws.Range("O23").Value = 0.01768353618551768
$ws.Range("P23").Value = 0.01768353618551768
$ws.Range("Q23").Value = 16.68796211465389
$ws.Range("R23").Value = 150.191659031885
$ws.Range("S23").Value = 0.007112198845974038
$ws.Range("T23").Value = 0.007112198845974036
$ws.Range("G24").Value = 23.12211433333333
$ws.Range("H24").Value = 69.366343
$ws.Range("I24").Value = 0.4021932475133977
$ws.Range("J24").Value = 0.4021932475133977
$ws.Range("M24").Value = 4.718667333333333
$ws.Range("N24").Value = 14.156002
$ws.Range("O24").Value = 0.1156146091272429
$ws.Range("P24").Value = 0.1156146091272429
$ws.Range("Q24").Value = 109.1055655822984
$ws.Range("R24").Value = 981.9500902406861
$ws.Range("S24").Value = 0.04649941510487792
$ws.Range("T24").Value = 0.04649941510487791
$ws.Range("G25").Value = 23.12211433333333
$ws.Range("H25").Value = 69.366343
$ws.Range("I25").Value = 0.4021932475133977
$ws.Range("J25").Value = 0.4021932475133977
$ws.Range("M25").Value = 27.97761966666667
$ws.Range("N25").Value = 83.93285900000001
$ws.Range("O25").Value = 0.6854947241613126
$ws.Range("P25").Value = 0.6854947241613125
$ws.Range("Q25").Value = 646.9017207071819
$ws.Range("R25").Value = 5822.115486364638
$ws.Range("S25").Value = 0.2757013492637391
$ws.Range("T25").Value = 0.275701349263739
$ws.Range("G26").Value = 23.12211433333333
$ws.Range("H26").Value = 69.366343
$ws.Range("I26").Value = 0.4021932475133977
$ws.Range("J26").Value = 0.4021932475133977
$ws.Range("M26").Value = 2.571656
$ws.Range("N26").Value = 7.714968
$ws.Range("O26").Value = 0.0630095283787885
$ws.Range("P26").Value = 0.0630095283787885
$ws.Range("Q26").Value = 59.46212405800267
$ws.Range("R26").Value = 535.159116522024
$ws.Range("S26").Value = 0.02534200684295254
$ws.Range("T26").Value = 0.02534200684295253
